# Updated parser to use TokenIteratorFieldRewriterSplit.
#
# The three Word fields in this document (the "{m: ... }" template
# instructions) are rewritten from real Word fields (w:fldChar / w:instrText)
# into plain literal text runs (w:t) that simply display the "{m: ...}"
# text - i.e. the field codes are "split" out into ordinary text, with the
# field delimiters removed and "{" / "}" added around the instruction text.

$d = $word.ActiveDocument

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-PackageXml([string]$bodyInner) {
    return "<pkg:package xmlns:pkg='http://schemas.microsoft.com/office/2006/xmlPackage'><pkg:part pkg:name='/word/document.xml' pkg:contentType='application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml'><pkg:xmlData><w:document $wNs><w:body>$bodyInner</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>"
}

$highlightRPr = "<w:rPr><w:color w:val=`"FFA500`"/><w:sz w:val=`"32`"/><w:highlight w:val=`"lightGray`"/></w:rPr>"

# --- Paragraph 2: "{m:let self = self.name}    <---    The variable mask an existing variable (self)."
$para2Inner = ""
$para2Inner += "<w:r><w:t>{m:</w:t></w:r>"
$para2Inner += "<w:r><w:t>let</w:t></w:r>"
$para2Inner += "<w:r><w:t xml:space=`"preserve`"> self</w:t></w:r>"
$para2Inner += "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
$para2Inner += "<w:r><w:t>=</w:t></w:r>"
$para2Inner += "<w:r><w:t xml:space=`"preserve`"> self.</w:t></w:r>"
$para2Inner += "<w:r><w:t>name}</w:t></w:r>"
$para2Inner += "<w:r><w:t xml:space=`"preserve`">    </w:t></w:r>"
$para2Inner += "<w:r>$highlightRPr<w:t>&lt;---</w:t></w:r>"
$para2Inner += "<w:r>$highlightRPr<w:t>The variable mask an existing variable (self).</w:t></w:r>"

$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertXML((New-PackageXml("<w:p>$para2Inner</w:p>")))

# --- Paragraph 3: "name = {m:self} ,"  (keep the "name = " text and trailing "," runs, rewrite field)
$para3Inner = ""
$para3Inner += "<w:proofErr w:type=`"spellStart`"/>"
$para3Inner += "<w:r><w:t>name</w:t></w:r>"
$para3Inner += "<w:proofErr w:type=`"spellEnd`"/>"
$para3Inner += "<w:r><w:t> </w:t></w:r>"
$para3Inner += "<w:proofErr w:type=`"gramStart`"/>"
$para3Inner += "<w:r><w:t>=</w:t></w:r>"
$para3Inner += "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>"
$para3Inner += "<w:proofErr w:type=`"gramEnd`"/>"
$para3Inner += "<w:r><w:t>{</w:t></w:r>"
$para3Inner += "<w:r><w:t>m</w:t></w:r>"
$para3Inner += "<w:r><w:t>:self</w:t></w:r>"
$para3Inner += "<w:r><w:t xml:space=`"preserve`">}</w:t></w:r>"
$para3Inner += "<w:r><w:t>,</w:t></w:r>"

$p3 = $d.Paragraphs.Item(3)
$p3.Range.InsertXML((New-PackageXml("<w:p>$para3Inner</w:p>")))

# --- Paragraph 4: "{m:endlet}"
$para4Inner = ""
$para4Inner += "<w:r><w:t>{</w:t></w:r>"
$para4Inner += "<w:r><w:t>m:</w:t></w:r>"
$para4Inner += "<w:r><w:t>endlet</w:t></w:r>"
$para4Inner += "<w:r><w:t xml:space=`"preserve`">}</w:t></w:r>"

$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertXML((New-PackageXml("<w:p>$para4Inner</w:p>")))

Write-Output "Done."
